# Update crypto price/volume data scraped on Wed Dec 13 11:56:14 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''41.091.30'
$ws.Range("E2").Value = '  -1.60%  '

# Row 3
$ws.Range("D3").Value = '''2.173.90'

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").Value = '''249.90'
$ws.Range("E5").Value = '  -0.53%  '

# Row 6
$ws.Range("E6").Value = '  -3.17%  '

# Row 7
$ws.Range("D7").Value = '''65.96'
$ws.Range("E7").Value = '  -7.51%  '

# Row 8
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("E9").Value = '  -2.17%  '

# Row 10
$ws.Range("D10").Value = '''58.64'
$ws.Range("E10").Value = '  +0.56%  '

# Row 11
$ws.Range("D11").Value = '''36.18'
$ws.Range("E11").Value = '  -10.85%  '

# Row 12
$ws.Range("E12").Value = '  -3.19%  '

# Row 13
$ws.Range("E13").Value = '  -1.46%  '

# Row 14
$ws.Range("D14").Value = '''6.83'
$ws.Range("E14").Value = '  -5.17%  '

# Row 15
$ws.Range("D15").Value = '''2.498.60'
$ws.Range("E15").Value = '  -1.85%  '

# Row 16
$ws.Range("D16").Value = '''14.18'
$ws.Range("E16").Value = '  -4.87%  '

# Row 17
$ws.Range("D17").Value = '''0.841'
$ws.Range("E17").Value = '  -3.60%  '

# Row 18
$ws.Range("D18").Value = '''2.180.51'
$ws.Range("E18").Value = '  -1.24%  '

# Row 19
$ws.Range("D19").Value = '''40.987.91'
$ws.Range("E19").Value = '  -1.65%  '

# Row 20
$ws.Range("E20").Value = '  -1.96%  '

# Row 21
$ws.Range("D21").Value = '''71.32'
$ws.Range("E21").Value = '  -1.85%  '

# Row 22
$ws.Range("D22").Value = '''6.02'
$ws.Range("E22").Value = '  -3.01%  '

# Row 23
$ws.Range("D23").Value = '''228.99'
$ws.Range("E23").Value = '  -2.35%  '

# Row 24
$ws.Range("E24").Value = '  -2.19%  '

# Row 25
$ws.Range("E25").Value = '  -4.94%  '

# Row 26
$ws.Range("E26").Value = '  +0.11%  '

# Row 27
$ws.Range("D27").Value = '''11.19'
$ws.Range("E27").Value = '  +1.38%  '

# Row 28
$ws.Range("E28").Value = '  -5.27%  '

# Row 29
$ws.Range("D29").Value = '''167.75'
$ws.Range("E29").Value = '  -1.49%  '

# Row 30
$ws.Range("E30").Value = '  -8.54%  '

# Row 31
$ws.Range("D31").Value = '''20.09'
$ws.Range("E31").Value = '  -3.02%  '

# Row 32
$ws.Range("E32").Value = '  -2.27%  '

# Row 33
$ws.Range("E33").Value = '  +1.75%  '

# Row 34
$ws.Range("D34").Value = '''0.0741'
$ws.Range("E34").Value = '  +0.51%  '

# Row 35
$ws.Range("E35").Value = '  -2.19%  '

# Row 36
$ws.Range("E36").Value = '  -4.92%  '

# Row 37
$ws.Range("D37").Value = '''3.95'
$ws.Range("E37").Value = '  -1.91%  '

# Row 38
$ws.Range("E38").Value = '  -8.04%  '

# Row 39
$ws.Range("D39").Value = '''0.0301'
$ws.Range("E39").Value = '  +0.83%  '

# Row 40
$ws.Range("D40").Value = '''5.45'
$ws.Range("E40").Value = '  +12.96%  '

# Row 41
$ws.Range("D41").Value = '''2.19'
$ws.Range("E41").Value = '  -3.85%  '

# Row 42
$ws.Range("D42").Value = '''5.49'
$ws.Range("E42").Value = '  -7.35%  '

# Row 43
$ws.Range("D43").Value = '''60.31'
$ws.Range("E43").Value = '  -8.09%  '

# Row 44
$ws.Range("D44").Value = '''11.14'
$ws.Range("E44").Value = '  -10.54%  '

# Row 45
$ws.Range("D45").Value = '''8.44'

# Row 46 (BinanceUSD)
$ws.Range("B46").Value = 'BinanceUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D46").Value = '''1.00'
$ws.Range("E46").Value = '  -0.16%  '

# Row 47 (Cronos)
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '''0.0988'
$ws.Range("E47").Value = '  -2.86%  '

# Row 48
$ws.Range("D48").Value = '''0.187'
$ws.Range("E48").Value = '  -8.57%  '

# Row 49
$ws.Range("E49").Value = '  -3.01%  '

# Row 50
$ws.Range("D50").Value = '''4.25'
$ws.Range("E50").Value = '  -10.33%  '

# Row 51
$ws.Range("E51").Value = '  -4.00%  '
